$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "devon_properties"

$ws.Cells.Item(3, 1).Value = "Rural Cabin Stunning Views"
$ws.Cells.Item(3, 2).Value = ""
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 4).Value = 197

$ws.Cells.Item(4, 1).Value = "The Owl's Nest"
$ws.Cells.Item(4, 2).Value = 4.99
$ws.Cells.Item(4, 3).Value = 143
$ws.Cells.Item(4, 4).Value = 314

$ws.Cells.Item(5, 1).Value = "A superb one bedroom apartment with sea views."
$ws.Cells.Item(5, 2).Value = 4.99
$ws.Cells.Item(5, 3).Value = 194
$ws.Cells.Item(5, 4).Value = 152

$ws.Cells.Item(6, 1).Value = "Countryside Cabin with Hot Tub and Tree Deck"
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(6, 3).Value = 356
$ws.Cells.Item(6, 4).Value = 165

$ws.Cells.Item(7, 1).Value = "No Snakes on this Plane!"
$ws.Cells.Item(7, 2).Value = 5
$ws.Cells.Item(7, 3).Value = 5
$ws.Cells.Item(7, 4).Value = 147

$ws.Cells.Item(8, 1).Value = "The Cabin at Axe View Hideout (hot tub stay)"
$ws.Cells.Item(8, 2).Value = 4.91
$ws.Cells.Item(8, 3).Value = 150
$ws.Cells.Item(8, 4).Value = 114

$ws.Cells.Item(9, 1).Value = "Luxury Roundhouse with log fired hot tub"
$ws.Cells.Item(9, 2).Value = 4.97
$ws.Cells.Item(9, 3).Value = 38
$ws.Cells.Item(9, 4).Value = 206

$ws.Cells.Item(10, 1).Value = "Hattie - luxury secluded coastal shepherds hut"
$ws.Cells.Item(10, 2).Value = 4.97
$ws.Cells.Item(10, 3).Value = 71
$ws.Cells.Item(10, 4).Value = 71

$ws.Cells.Item(11, 1).Value = "''Rockpool' is a 15 minute walk to Bantham Beach."
$ws.Cells.Item(11, 2).Value = 4.99
$ws.Cells.Item(11, 3).Value = 300
$ws.Cells.Item(11, 4).Value = 69

$ws.Cells.Item(12, 1).Value = "North Devon: Treetops - Peace in Nature"
$ws.Cells.Item(12, 2).Value = 4.95
$ws.Cells.Item(12, 3).Value = 83
$ws.Cells.Item(12, 4).Value = 91

$ws.Cells.Item(13, 1).Value = "Tranquil Room by the Green Field"
$ws.Cells.Item(13, 2).Value = 5
$ws.Cells.Item(13, 3).Value = 10
$ws.Cells.Item(13, 4).Value = 59

$ws.Cells.Item(14, 1).Value = "The Posh Shed"
$ws.Cells.Item(14, 2).Value = 4.93
$ws.Cells.Item(14, 3).Value = 252
$ws.Cells.Item(14, 4).Value = 74

$ws.Cells.Item(15, 1).Value = "The Shippon. Unique luxurious South Devon getaway."
$ws.Cells.Item(15, 2).Value = 4.96
$ws.Cells.Item(15, 3).Value = 330
$ws.Cells.Item(15, 4).Value = 216

$ws.Cells.Item(16, 1).Value = "Little Owl birdhouse: Stunning sea views & beach"
$ws.Cells.Item(16, 2).Value = 4.99
$ws.Cells.Item(16, 3).Value = 73
$ws.Cells.Item(16, 4).Value = 182

$ws.Cells.Item(17, 1).Value = "Ashmead Shepherds Hut"
$ws.Cells.Item(17, 2).Value = 5
$ws.Cells.Item(17, 3).Value = 11
$ws.Cells.Item(17, 4).Value = 128

$ws.Cells.Item(18, 1).Value = "Wolf Valley- 'The Coracle' geodesic dome ~pondside"
$ws.Cells.Item(18, 2).Value = 4.95
$ws.Cells.Item(18, 3).Value = 175
$ws.Cells.Item(18, 4).Value = 140

$ws.Cells.Item(19, 1).Value = "Luxury Shepherds Hut with wood fired hot tub"
$ws.Cells.Item(19, 2).Value = 5
$ws.Cells.Item(19, 3).Value = 129
$ws.Cells.Item(19, 4).Value = 161

$ws.Cells.Item(20, 1).Value = "The Lodge"
$ws.Cells.Item(20, 2).Value = 5
$ws.Cells.Item(20, 3).Value = 80
$ws.Cells.Item(20, 4).Value = 154

$ws.Cells.Item(21, 1).Value = "Lilypod Heron –Luxury Floating Dome Stay in Devon"
$ws.Cells.Item(21, 2).Value = 5
$ws.Cells.Item(21, 3).Value = 30
$ws.Cells.Item(21, 4).Value = 258

$ws.Cells.Item(22, 1).Value = "The Cabin Devon rural retreat perfect for couples."
$ws.Cells.Item(22, 2).Value = 5
$ws.Cells.Item(22, 3).Value = 409
$ws.Cells.Item(22, 4).Value = 112

$ws.Cells.Item(23, 1).Value = "Cosy traditional Devon cottage"
$ws.Cells.Item(23, 2).Value = 5
$ws.Cells.Item(23, 3).Value = 6
$ws.Cells.Item(23, 4).Value = 64

$ws.Cells.Item(24, 1).Value = "''The Weekender' @Cleavefarmcottages, Crackington"
$ws.Cells.Item(24, 2).Value = 4.99
$ws.Cells.Item(24, 3).Value = 511
$ws.Cells.Item(24, 4).Value = 100

$ws.Cells.Item(25, 1).Value = "Cornwall Beach House - Panoramic Sea Views"
$ws.Cells.Item(25, 2).Value = 4.89
$ws.Cells.Item(25, 3).Value = 56
$ws.Cells.Item(25, 4).Value = 419

$ws.Cells.Item(26, 1).Value = "Glamping retreat: dome & wagon with alpacas, Devon"
$ws.Cells.Item(26, 2).Value = 5
$ws.Cells.Item(26, 3).Value = 64
$ws.Cells.Item(26, 4).Value = 182

$ws.Cells.Item(27, 1).Value = "The Maple Room, Totnes, Guest Suite own entrance."
$ws.Cells.Item(27, 2).Value = 4.91
$ws.Cells.Item(27, 3).Value = 451
$ws.Cells.Item(27, 4).Value = 40

$ws.Cells.Item(28, 1).Value = "Finest Retreats | Yeworthy Eco-Treehouse"
$ws.Cells.Item(28, 2).Value = 4.94
$ws.Cells.Item(28, 3).Value = 390
$ws.Cells.Item(28, 4).Value = 152

$ws.Cells.Item(29, 1).Value = "Idyllic Secluded Pondside Cabin-Devon Countryside"
$ws.Cells.Item(29, 2).Value = 4.99
$ws.Cells.Item(29, 3).Value = 113
$ws.Cells.Item(29, 4).Value = 111

$ws.Cells.Item(30, 1).Value = "Swallow View, Umberleigh, North Devon"
$ws.Cells.Item(30, 2).Value = 4.99
$ws.Cells.Item(30, 3).Value = 281
$ws.Cells.Item(30, 4).Value = 75

$ws.Cells.Item(31, 1).Value = "Romantic Ocean View Couples Retreat Cornwall"
$ws.Cells.Item(31, 2).Value = 4.99
$ws.Cells.Item(31, 3).Value = 150
$ws.Cells.Item(31, 4).Value = 357

$ws.Cells.Item(32, 1).Value = "Cosy Shepherd's Hut in beautiful North Devon"
$ws.Cells.Item(32, 2).Value = 4.97
$ws.Cells.Item(32, 3).Value = 145
$ws.Cells.Item(32, 4).Value = 95

$ws.Cells.Item(33, 1).Value = "Idyllic Shepherd Hut in Dartmoor"
$ws.Cells.Item(33, 2).Value = 4.96
$ws.Cells.Item(33, 3).Value = 78
$ws.Cells.Item(33, 4).Value = 101

$ws.Cells.Item(34, 1).Value = "The Wizards Cauldron -Harry Potter Themed"
$ws.Cells.Item(34, 2).Value = 5
$ws.Cells.Item(34, 3).Value = 163
$ws.Cells.Item(34, 4).Value = 176

$ws.Cells.Item(35, 1).Value = "Coastal cliff top chalet, within a private garden"
$ws.Cells.Item(35, 2).Value = 4.97
$ws.Cells.Item(35, 3).Value = 30
$ws.Cells.Item(35, 4).Value = 135

$ws.Cells.Item(36, 1).Value = "The Old Stables - A Cosy Riverside Retreat"
$ws.Cells.Item(36, 2).Value = 4.98
$ws.Cells.Item(36, 3).Value = 206
$ws.Cells.Item(36, 4).Value = 107

$ws.Cells.Item(37, 1).Value = "Ashridge Farm"
$ws.Cells.Item(37, 2).Value = 4.84
$ws.Cells.Item(37, 3).Value = 94
$ws.Cells.Item(37, 4).Value = 51

$ws.Cells.Item(38, 1).Value = "The Granary"
$ws.Cells.Item(38, 2).Value = 4.98
$ws.Cells.Item(38, 3).Value = 115
$ws.Cells.Item(38, 4).Value = 110

$ws.Cells.Item(39, 1).Value = "Kingfisher Pod: Scenic Glamping at Milemead Lakes"
$ws.Cells.Item(39, 2).Value = 4.93
$ws.Cells.Item(39, 3).Value = 389
$ws.Cells.Item(39, 4).Value = 64

$ws.Cells.Item(40, 1).Value = "Fantastic coast and country retreat."
$ws.Cells.Item(40, 2).Value = 4.93
$ws.Cells.Item(40, 3).Value = 121
$ws.Cells.Item(40, 4).Value = 111

$ws.Cells.Item(41, 1).Value = "Tranquil Shepherd's Hut with hot tub access [DWK]"
$ws.Cells.Item(41, 2).Value = 4.91
$ws.Cells.Item(41, 3).Value = 126
$ws.Cells.Item(41, 4).Value = 78

$ws.Cells.Item(42, 1).Value = "Honeysuckle Shepherd Hut~Secluded ~Luxury~Hot Tub"
$ws.Cells.Item(42, 2).Value = 5
$ws.Cells.Item(42, 3).Value = 237
$ws.Cells.Item(42, 4).Value = 152

$ws.Cells.Item(43, 1).Value = "Treetop cabin & outdoor bath in 45 acre woodland"
$ws.Cells.Item(43, 2).Value = 5
$ws.Cells.Item(43, 3).Value = 32
$ws.Cells.Item(43, 4).Value = 246

$ws.Cells.Item(44, 1).Value = "‘The Old Laundry Room’ Unique Space"
$ws.Cells.Item(44, 2).Value = 4.91
$ws.Cells.Item(44, 3).Value = 258
$ws.Cells.Item(44, 4).Value = 99

$ws.Cells.Item(45, 1).Value = "Cosy barn with hot tub and alpacas"
$ws.Cells.Item(45, 2).Value = 5
$ws.Cells.Item(45, 3).Value = 59
$ws.Cells.Item(45, 4).Value = 164

$ws.Cells.Item(46, 1).Value = "Luxury Cabin Retreat with Hot Tub - Langman"
$ws.Cells.Item(46, 2).Value = 4.99
$ws.Cells.Item(46, 3).Value = 322
$ws.Cells.Item(46, 4).Value = 203

$ws.Cells.Item(47, 1).Value = "Handcrafted hut with outdoor bath"
$ws.Cells.Item(47, 2).Value = 5
$ws.Cells.Item(47, 3).Value = 33
$ws.Cells.Item(47, 4).Value = 183

$ws.Cells.Item(48, 1).Value = "The Drey Near Braunton NorthDevon romantic retreat"
$ws.Cells.Item(48, 2).Value = 5
$ws.Cells.Item(48, 3).Value = 261
$ws.Cells.Item(48, 4).Value = 188

$ws.Cells.Item(49, 1).Value = "The Garden Studio at the Tithe Barn"
$ws.Cells.Item(49, 2).Value = 4.8
$ws.Cells.Item(49, 3).Value = 343
$ws.Cells.Item(49, 4).Value = 72

$ws.Cells.Item(50, 1).Value = "Shepherds hut hot tub & Firepit all wood included"
$ws.Cells.Item(50, 2).Value = 4.98
$ws.Cells.Item(50, 3).Value = 176
$ws.Cells.Item(50, 4).Value = 119

$ws.Cells.Item(51, 1).Value = "Valley View tranquillity nr Pigs Nose"
$ws.Cells.Item(51, 2).Value = 4.98
$ws.Cells.Item(51, 3).Value = 123
$ws.Cells.Item(51, 4).Value = 182

$ws.Cells.Item(52, 1).Value = "1950's bungalow with lake views"
$ws.Cells.Item(52, 2).Value = 5
$ws.Cells.Item(52, 3).Value = 3
$ws.Cells.Item(52, 4).Value = 127
